{"js": "// Adds three new paragraphs to the end of the document body, after the\n// paragraph that currently ends with \"...Thank you, Nadiya!\":\n//   1. An empty paragraph (blank line separator)\n//   2. A paragraph honoring the Natural History Museum volunteers\n//   3. A paragraph finishing that thought (\"Year! Thank you ladies...\")\n//\n// All new paragraphs use the same paragraph style (\"Normal\" / style0) as\n// the rest of the document's body paragraphs.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nlastParagraph.load(\"style\");\nawait context.sync();\n\nconst baseStyle = lastParagraph.style;\n\n// 1) Blank separator paragraph.\nconst blankParagraph = lastParagraph.insertParagraph(\"\", \"After\");\nblankParagraph.style = baseStyle;\n\n// 2) Paragraph introducing the Natural History Museum \"outreach ladies\".\nconst para2Text =\n  \"Additionally, the Hart Honeree  for the Natural History Museum's 2013 \" +\n  \"Volunteer of the Year award goes to a group of volunteers that help \" +\n  \"out with our program to provide an immersive education experience for \" +\n  \"local elementary schools.  We are so delighted to honor Judy Beadles, \" +\n  \"Phyllis Dozier, Janet Labick, Diane Whetzel, and Joy Wilson \\u2014 the \" +\n  \"\\u201coutreach ladies\\u201d as the NHM Volunteer of the \";\nconst para2 = blankParagraph.insertParagraph(para2Text, \"After\");\npara2.style = baseStyle;\n\n// 3) Closing paragraph.\nconst para3Text = \"Year! Thank you ladies for everything you do!\";\nconst para3 = para2.insertParagraph(para3Text, \"After\");\npara3.style = baseStyle;\n\nawait context.sync();\n", "ps1": "# Adds three new paragraphs to the end of the document body, after the\n# paragraph that currently ends with \"...Thank you, Nadiya!\":\n#   1. An empty paragraph (blank line separator)\n#   2. A paragraph honoring the Natural History Museum volunteers\n#   3. A paragraph finishing that thought (\"Year! Thank you ladies...\")\n#\n# All new paragraphs keep the same paragraph style (\"Normal\" / style0) that\n# Word automatically carries forward from the preceding paragraph.\n\n$d = $word.ActiveDocument\n\n$emDash = [char]0x2014\n$leftQuote = [char]0x201C\n$rightQuote = [char]0x201D\n\n$para2Text = \"Additionally, the Hart Honeree  for the Natural History Museum's 2013 Volunteer of the Year award goes to a group of volunteers that help out with our program to provide an immersive education experience for local elementary schools.  We are so delighted to honor Judy Beadles, Phyllis Dozier, Janet Labick, Diane Whetzel, and Joy Wilson \" + $emDash + \" the \" + $leftQuote + \"outreach ladies\" + $rightQuote + \" as the NHM Volunteer of the \"\n$para3Text = \"Year! Thank you ladies for everything you do!\"\n\n# 1) Blank separator paragraph.\n$rng = $d.Paragraphs.Last.Range\n$rng.Collapse(0)\n$rng.InsertParagraphAfter()\n\n# 2) Paragraph introducing the Natural History Museum \"outreach ladies\".\n$rng = $d.Paragraphs.Last.Range\n$rng.Collapse(0)\n$rng.InsertParagraphAfter()\n\n$rng = $d.Paragraphs.Last.Range\n$rng.Collapse(0)\n$rng.InsertAfter($para2Text)\n\n# 3) Closing paragraph.\n$rng = $d.Paragraphs.Last.Range\n$rng.Collapse(0)\n$rng.InsertParagraphAfter()\n\n$rng = $d.Paragraphs.Last.Range\n$rng.Collapse(0)\n$rng.InsertAfter($para3Text)\n"}
